$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3760.6924
$ws.Range("I2").Value = 6142
$ws.Range("J2").Value = 982.5
$ws.Range("K2").Value = 6142
$ws.Range("L2").Value = 982.5
$ws.Range("M2").Value = -6029
$ws.Range("N2").Value = -1208.5

$ws.Range("H41").Value = 1236.7778
$ws.Range("I41").Value = 1126
$ws.Range("J41").Value = 1268.4286
$ws.Range("K41").Value = 1126
$ws.Range("L41").Value = 1268.4286
$ws.Range("M41").Value = -686
$ws.Range("N41").Value = -2148.4286

$ws.Range("H53").Value = 4810.9546
$ws.Range("I53").Value = 263.46667
$ws.Range("J53").Value = 14555.571
$ws.Range("K53").Value = 263.46667
$ws.Range("L53").Value = 14555.571
$ws.Range("M53").Value = 373.53333
$ws.Range("N53").Value = -15829.571

$ws.Range("H64").Value = 8961.462
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8961.462
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 8961.462
$ws.Range("N64").Value = -9457.462

$ws.Range("H67").Value = 8961.462
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8961.462
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 8961.462
$ws.Range("N67").Value = -10677.462

$ws.Range("H74").Value = 6957.4443
$ws.Range("I74").Value = 4574.5
$ws.Range("J74").Value = 7638.2856
$ws.Range("K74").Value = 4574.5
$ws.Range("L74").Value = 7638.2856
$ws.Range("M74").Value = -3638.5
$ws.Range("N74").Value = -9510.285599999999

$ws.Range("H76").Value = 5604.2666
$ws.Range("I76").Value = 4569
$ws.Range("J76").Value = 5980.727
$ws.Range("K76").Value = 4569
$ws.Range("L76").Value = 5980.727
$ws.Range("M76").Value = -4254
$ws.Range("N76").Value = -6610.727

$ws.Range("H77").Value = 6957.4443
$ws.Range("I77").Value = 4574.5
$ws.Range("J77").Value = 7638.2856
$ws.Range("K77").Value = 22872.5
$ws.Range("L77").Value = 38191.428
$ws.Range("M77").Value = -18192.5
$ws.Range("N77").Value = -47551.428

$ws.Range("H79").Value = 5604.2666
$ws.Range("I79").Value = 4569
$ws.Range("J79").Value = 5980.727
$ws.Range("K79").Value = 4569
$ws.Range("L79").Value = 5980.727
$ws.Range("M79").Value = -3477
$ws.Range("N79").Value = -8164.727

$ws.Range("H86").Value = 5070.36
$ws.Range("I86").Value = 3773.9167
$ws.Range("J86").Value = 6267.077
$ws.Range("K86").Value = 3773.9167
$ws.Range("L86").Value = 6267.077
$ws.Range("M86").Value = -2650.9167
$ws.Range("N86").Value = -8513.077000000001

$ws.Range("H89").Value = 5070.36
$ws.Range("I89").Value = 3773.9167
$ws.Range("J89").Value = 6267.077
$ws.Range("K89").Value = 18869.5835
$ws.Range("L89").Value = 31335.385
$ws.Range("M89").Value = -13253.5835
$ws.Range("N89").Value = -42567.385

$ws.Range("H92").Value = 1333
$ws.Range("I92").Value = 401.5909
$ws.Range("J92").Value = 4748.1665
$ws.Range("K92").Value = 401.5909
$ws.Range("L92").Value = 4748.1665
$ws.Range("M92").Value = 846.4091000000001
$ws.Range("N92").Value = -7244.1665

$ws.Range("H107").Value = 23881856
$ws.Range("I107").Value = 37038250
$ws.Range("J107").Value = 200350.4
$ws.Range("K107").Value = 37038250
$ws.Range("L107").Value = 200350.4
$ws.Range("M107").Value = -37036330
$ws.Range("N107").Value = -204190.4

$ws.Range("H132").Value = 15626871
$ws.Range("I132").Value = 16668635
$ws.Range("J132").Value = 411.5
$ws.Range("K132").Value = 50005905
$ws.Range("L132").Value = 1234.5
$ws.Range("M132").Value = -50003375
$ws.Range("N132").Value = -6294.5

$ws.Range("H137").Value = 37361.84
$ws.Range("I137").Value = 44824.145
$ws.Range("J137").Value = 3366.889
$ws.Range("K137").Value = 134472.435
$ws.Range("L137").Value = 10100.667
$ws.Range("M137").Value = -131922.435
$ws.Range("N137").Value = -15200.667


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13144.682
$ws.Range("I32").Value = 8272.532999999999
$ws.Range("J32").Value = 23585
$ws.Range("K32").Value = 8272.532999999999
$ws.Range("L32").Value = 23585
$ws.Range("M32").Value = -7985.532999999999
$ws.Range("N32").Value = -24159


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3226876.8
$ws.Range("I86").Value = 5000971
$ws.Range("J86").Value = 1251.1818
$ws.Range("K86").Value = 5000971
$ws.Range("L86").Value = 1251.1818
$ws.Range("M86").Value = -4999848
$ws.Range("N86").Value = -3497.1818

$ws.Range("H89").Value = 3226876.8
$ws.Range("I89").Value = 5000971
$ws.Range("J89").Value = 1251.1818
$ws.Range("K89").Value = 25004855
$ws.Range("L89").Value = 6255.909000000001
$ws.Range("M89").Value = -24999239
$ws.Range("N89").Value = -17487.909

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws.Range("N126").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23526.809
$ws.Range("I31").Value = 2512.6875
$ws.Range("J31").Value = 34372.805
$ws.Range("K31").Value = 2512.6875
$ws.Range("L31").Value = 34372.805
$ws.Range("M31").Value = -2217.6875
$ws.Range("N31").Value = -34962.805

$ws.Range("H34").Value = 23526.809
$ws.Range("I34").Value = 2512.6875
$ws.Range("J34").Value = 34372.805
$ws.Range("K34").Value = 2512.6875
$ws.Range("L34").Value = 34372.805
$ws.Range("M34").Value = -2310.6875
$ws.Range("N34").Value = -34776.805

$ws.Range("H58").Value = 6854.75
$ws.Range("I58").Value = 9228.857
$ws.Range("J58").Value = 3531
$ws.Range("K58").Value = 9228.857
$ws.Range("L58").Value = 3531
$ws.Range("M58").Value = -9025.857
$ws.Range("N58").Value = -3937

$ws.Range("H92").Value = 37516.5
$ws.Range("I92").Value = 25000
$ws.Range("J92").Value = 40019.8
$ws.Range("K92").Value = 25000
$ws.Range("L92").Value = 40019.8
$ws.Range("M92").Value = -22504
$ws.Range("N92").Value = -45011.8

$ws.Range("H136").Value = 6854.75
$ws.Range("I136").Value = 9228.857
$ws.Range("J136").Value = 3531
$ws.Range("K136").Value = 27686.571
$ws.Range("L136").Value = 10593
$ws.Range("M136").Value = -25136.571
$ws.Range("N136").Value = -15693


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5287.095
$ws.Range("I81").Value = 617.8
$ws.Range("J81").Value = 6746.25
$ws.Range("K81").Value = 1853.4
$ws.Range("L81").Value = 20238.75
$ws.Range("M81").Value = -730.3999999999999
$ws.Range("N81").Value = -22484.75

$ws.Range("H84").Value = 5287.095
$ws.Range("I84").Value = 617.8
$ws.Range("J84").Value = 6746.25
$ws.Range("K84").Value = 5560.2
$ws.Range("L84").Value = 60716.25
$ws.Range("M84").Value = 55.80000000000018
$ws.Range("N84").Value = -71948.25

$ws.Range("H122").Value = 1460.125
$ws.Range("I122").Value = 1382.2858
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 12440.5722
$ws.Range("L122").Value = 18045
$ws.Range("M122").Value = -9990.572200000001
$ws.Range("N122").Value = -22945

$ws.Range("H128").Value = 199989.5
$ws.Range("I128").Value = 199989.5
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 599968.5
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -594988.5

$ws.Range("H132").Value = 1818.091
$ws.Range("I132").Value = 1599.5
$ws.Range("J132").Value = 1866.6666
$ws.Range("K132").Value = 14395.5
$ws.Range("L132").Value = 16799.9994
$ws.Range("M132").Value = -11865.5
$ws.Range("N132").Value = -21859.9994

$ws.Range("H141").Value = 2949.5
$ws.Range("I141").Value = 2949.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8848.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3668.5


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3826.4736
$ws.Range("I100").Value = 3307.923
$ws.Range("J100").Value = 4950
$ws.Range("K100").Value = 3307.923
$ws.Range("L100").Value = 4950
$ws.Range("M100").Value = -2766.923
$ws.Range("N100").Value = -6032


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3026.0908
$ws.Range("I122").Value = 2188.6843
$ws.Range("J122").Value = 8329.666999999999
$ws.Range("K122").Value = 6566.0529
$ws.Range("L122").Value = 24989.001
$ws.Range("M122").Value = -4116.0529
$ws.Range("N122").Value = -29889.001

$ws.Range("H137").Value = 69999.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 69999.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 69999.5
$ws.Range("N137").Value = -80199.5

